$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 21: the "Time Spent (Hours)" entry for day 18 (20/3/2024) was
# overstated - correct it from 40 down to 4.
$ws.Range("C21").Value = 4

# New log entry - day 19, dated 21/3/2024, 8 hours spent, describing the
# newly added product search engine functionality.
$ws.Cells.Item(22, 1).Value = 19
$ws.Cells.Item(22, 2).Value = "21/3/2024"
$ws.Cells.Item(22, 3).Value = 8
$ws.Cells.Item(22, 4).Value = "Added search engine functionality for products"

# Match the rest of the table's centered alignment (style used by every
# other data row) for the freshly added row.
$newRow = $ws.Range("A22:D22")
$newRow.HorizontalAlignment = -4108  # xlCenter
$newRow.VerticalAlignment = -4108    # xlCenter

# Scroll the view down a bit and move the active selection to reflect
# where the user was working after adding the new row.
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C25").Select()
